$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B3" = 81716.10000000001
    "C3" = 76255.21000000001
    "D3" = 76187.89999999999
    "E3" = 68580.99000000001
    "F3" = 110584.34
    "G3" = 413324.54

    "B4" = -23546.99
    "C4" = -12433.3
    "D4" = -16701.58
    "E4" = -17088.3
    "F4" = -19492.48
    "G4" = -89262.64999999999

    "B5" = -24322.07
    "C5" = -17363.31
    "D5" = -15047.77
    "E5" = -18153.58
    "F5" = -18839.59
    "G5" = -93726.32000000001

    "B6" = -12972.63
    "C6" = -11891.43
    "D6" = -15180.02
    "E6" = -14769.04
    "F6" = -12670.83
    "G6" = -67483.95

    "B7" = -11891.18
    "C7" = -20277.56
    "D7" = -16132.08
    "E7" = -18620.19
    "F7" = -13311.51
    "G7" = -80232.52

    "B8" = 57327.05
    "C8" = 69089.91
    "D8" = 46505.97
    "E8" = 68042.99000000001
    "F8" = 97647.47
    "G8" = 338613.39

    "B9" = -22686.33
    "C9" = -12307.34
    "D9" = -20829.45
    "E9" = -15279.92
    "F9" = -17002.62
    "G9" = -88105.66

    "B10" = -18135.67
    "C10" = -16482.58
    "D10" = -15629.56
    "E10" = -19907.67
    "F10" = -12846.96
    "G10" = -83002.44

    "B11" = -24328.03
    "C11" = -23259.89
    "D11" = -18180.72
    "E11" = -13566.12
    "F11" = -12676.5
    "G11" = -92011.25999999999

    "B12" = -13666.72
    "C12" = -18150.04
    "D12" = -17989.96
    "E12" = -23666.17
    "F12" = -17427.45
    "G12" = -90900.34

    "B13" = -16484.8
    "C13" = -14390.87
    "D13" = -18230.71
    "E13" = -30849.23
    "F13" = -13264.47
    "G13" = -93220.08

    "B14" = -11866.61
    "C14" = -20344.57
    "D14" = -15438.74
    "E14" = -15008.15
    "F14" = -18061.25
    "G14" = -80719.32000000001

    "B15" = -15592
    "C15" = -18102.8
    "D15" = -12464.41
    "E15" = -15161.04
    "F15" = -12295.02
    "G15" = -73615.27

    "B16" = -17838.35
    "C16" = -22945.35
    "D16" = -16119.28
    "E16" = -19638.93
    "F16" = -20240.81
    "G16" = -96782.72

    "B17" = -17905.11
    "C17" = -17650.67
    "D17" = -19095.81
    "E17" = -25506.55
    "F17" = -17105.98
    "G17" = -97264.12

    "B18" = 452802.19
    "C18" = 439948.04
    "D18" = 521104.54
    "E18" = 400785.61
    "F18" = 439761.29
    "G18" = 2254401.67

    "B19" = -14587.84
    "C19" = -19625.88
    "D19" = -20112.34
    "E19" = -18387.95
    "F19" = -16297.45
    "G19" = -89011.46000000001

    "B20" = 346021.01
    "C20" = 340067.57
    "D20" = 406645.98
    "E20" = 271806.75
    "F20" = 426460.18
    "G20" = 1791001.49
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
